$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "23.238.27"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.72%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.612.92"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.82%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.26%  "

$ws.Range("E5").Value = "  +0.11%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "302.70"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.58%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3783"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.11%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "51.68"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.56%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3529"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.78%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08089"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.13%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.205"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.92%  "

$ws.Range("E12").Value = "  +0.26%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.14"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.45%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.360"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.78%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.270"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.85%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001211"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.77%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.586.40"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.22%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "93.79"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.32%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06897"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.13%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.470"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.92%  "

$ws.Range("B21").Value = "Avalanche"
$ws.Range("C21").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.21"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.80%  "

$ws.Range("B22").Value = "Dai"
$ws.Range("C22").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.16%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.31"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.25%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "23.222.09"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.74%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.534"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.65%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.030"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -6.34%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.85"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.33%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "151.02"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.27%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.235"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.00%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "132.06"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.59%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.778.60"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.48%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.068"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +12.28%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.452"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.96%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.096"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -9.04%  "

$ws.Range("E35").Value = "  +2.95%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02704"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.20%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.08682"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.76%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2451"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.06%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06915"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.70%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.849"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.22%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.324"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.48%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6868"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.79%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "11.91"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.18%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.22"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.44%  "

$ws.Range("E45").Value = "  +0.09%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6303"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.43%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.939"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.46%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.248"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.26%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07868"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.54%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "127.60"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.54%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.168"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.85%  "
